$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-20 Monday" "2025-10-21 Tuesday"

Replace-Text "87×11=957" "17×72=1224"
Replace-Text "40×27=1080" "80×99=7920"
Replace-Text "25×63=1575" "58×98=5684"
Replace-Text "59×68=4012" "91×27=2457"
Replace-Text "95×84=7980" "40×45=1800"

Replace-Text "61×14=854" "90×30=2700"
Replace-Text "30×43=1290" "35×73=2555"
Replace-Text "78×61=4758" "37×36=1332"
Replace-Text "48×75=3600" "82×96=7872"
Replace-Text "90×56=5040" "97×17=1649"

Replace-Text "47×28=1316" "33×36=1188"
Replace-Text "38×14=532" "42×93=3906"
Replace-Text "56×27=1512" "69×52=3588"
Replace-Text "53×92=4876" "96×95=9120"
Replace-Text "71×63=4473" "90×32=2880"

Replace-Text "47×86=4042" "74×30=2220"
Replace-Text "47×85=3995" "17×14=238"
Replace-Text "28×47=1316" "27×99=2673"
Replace-Text "31×84=2604" "39×92=3588"
Replace-Text "21×27=567" "31×47=1457"

Replace-Text "76×37=2812" "90×13=1170"
Replace-Text "26×35=910" "39×70=2730"
Replace-Text "17×60=1020" "53×28=1484"
Replace-Text "75×49=3675" "97×84=8148"
Replace-Text "26×53=1378" "47×45=2115"
